$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-01 Thursday" "2024-02-02 Friday"

Replace-Text "28×95=2660" "56×70=3920"
Replace-Text "30×91=2730" "58×80=4640"
Replace-Text "69×60=4140" "34×82=2788"
Replace-Text "87×41=3567" "19×84=1596"
Replace-Text "92×98=9016" "27×27=729"

Replace-Text "48×78=3744" "40×91=3640"
Replace-Text "24×80=1920" "89×37=3293"
Replace-Text "49×89=4361" "36×64=2304"
Replace-Text "78×20=1560" "44×56=2464"
Replace-Text "87×47=4089" "23×90=2070"

Replace-Text "25×38=950" "32×98=3136"
Replace-Text "26×33=858" "86×60=5160"
Replace-Text "86×48=4128" "73×69=5037"
Replace-Text "47×58=2726" "44×31=1364"
Replace-Text "59×11=649" "98×15=1470"

Replace-Text "21×62=1302" "56×95=5320"
Replace-Text "68×63=4284" "44×84=3696"
Replace-Text "66×54=3564" "68×83=5644"
Replace-Text "55×67=3685" "79×40=3160"
Replace-Text "14×69=966" "35×55=1925"

Replace-Text "58×81=4698" "99×80=7920"
Replace-Text "38×86=3268" "34×30=1020"
Replace-Text "27×91=2457" "47×23=1081"
Replace-Text "99×94=9306" "98×13=1274"
Replace-Text "18×67=1206" "61×37=2257"
